$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Bad Drivers table updates
$ws.Range("C4").Value = 158
$ws.Range("D4").Value = 98.2
$ws.Range("C5").Value = 237

# Good Drivers table updates
$ws.Range("B13").Value = 449371
$ws.Range("B14").Value = 77999
